# Proyectos Informacionales.xlsx — add missing project row
# "Excels corregidos espacios quitados en proyectos de tls 201812 y 201811"
#
# A new project (EXT-000193-00894 / "DNA AM:POC MICRO FY '18'" / "AM") was
# missing from the table and gets inserted as row 14, pushing the existing
# rows 14-27 down to 15-28. The new row's Project cell (column A) is
# highlighted with a distinct font (Times New Roman 12pt, vertically
# centered) while the rest of the row keeps the table's normal formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14; this shifts rows 14:27 down to 15:28
# and extends the sheet dimension from D27 to D28 automatically.
[void]$ws.Rows("14:14").Insert()

# The new row should look like every other data row (same fill/border),
# so copy that formatting from B15:D15 (the row that used to be row 14)
# into the freshly inserted B14:D14.
[void]$ws.Range("B15:D15").Copy()
[void]$ws.Range("B14:D14").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new project's data. D14 (Area) is intentionally left blank.
$ws.Range("A14").Value2 = "EXT-000193-00894"
$ws.Range("B14").Value2 = "DNA AM:POC MICRO FY '18'"
$ws.Range("C14").Value2 = "AM"

# Give the new project's code cell its own highlight font.
$fnt = $ws.Range("A14").Font
$fnt.Size = 12
$fnt.Name = "Times New Roman"
$ws.Range("A14").VerticalAlignment = -4108 # xlVAlignCenter

# This row is a bit taller than the default to fit the larger font.
$ws.Rows("14:14").RowHeight = 15.5

# Leave the selection on the newly added row, like in the saved workbook.
[void]$ws.Range("B14").Select()

Write-Host "Inserted new project row EXT-000193-00894 at row 14"
